$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new BOM row above row 9 ("Input Gear") to make room for the
#    new "Motor Controller" line item. This pushes the old rows 9/10 (Input
#    Gear / Output Gear) down to 10/11 and the trailing SUM row from 20 to 21,
#    and Excel auto-extends the shared SUM range from D2:D10 to D2:D11.
# ---------------------------------------------------------------------------
$ws.Rows(9).Insert()

# ---------------------------------------------------------------------------
# 2. Motor (row 8): the actual motor purchased has the same unit price but a
#    different measured/actual total cost, and a different (correct) vendor
#    link (Crouzet motor on Digikey instead of the old Trinamic link).
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = 180.98

$ws.Range("E2").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = "https://www.digikey.ca/en/products/detail/crouzet/82890001/3190319"

# ---------------------------------------------------------------------------
# 3. New row 9: Motor Controller.
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "Motor Controller"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 29.95
$ws.Range("D9").Formula = "=C9"

$ws.Range("E2").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = "https://www.pololu.com/product/2991"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Move the picture's anchor down by one row, matching the row that was
#    just inserted above it (from row 11/51 to row 12/52, 0-indexed).
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = $shp.Top + 15

# ---------------------------------------------------------------------------
# 5. Selection lands on the new grand-total cell.
# ---------------------------------------------------------------------------
[void]$ws.Range("D21").Select()

$excel.CalculateFull()
